# Applies the scheduled market-data refresh to the Leve profit sheets.
# Columns H-N on each sheet are refreshed average-price / profit figures;
# a handful of rows gain or lose an HQ-profit (M) / HQ-loss (N) cell entirely
# depending on whether an HQ price was reported this run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1721.75
$ws.Range("I12").Value = 296.5
$ws.Range("K12").Value = 296.5
$ws.Range("M12").Value = -126.5
$ws.Range("H15").Value = 1465.5625
$ws.Range("I15").Value = 1465.5625
$ws.Range("K15").Value = 4396.6875
$ws.Range("M15").Value = -4227.6875
$ws.Range("H51").Value = 4500
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H64").Value = 4866.3335
$ws.Range("J64").Value = 5214.143
$ws.Range("L64").Value = 5214.143
$ws.Range("N64").Value = -5710.143
$ws.Range("H67").Value = 4866.3335
$ws.Range("J67").Value = 5214.143
$ws.Range("L67").Value = 5214.143
$ws.Range("N67").Value = -6930.143
$ws.Range("H96").Value = 14503
$ws.Range("I96").Value = 26756.5
$ws.Range("J96").Value = 2249.5
$ws.Range("K96").Value = 80269.5
$ws.Range("L96").Value = 6748.5
$ws.Range("M96").Value = -78896.5
$ws.Range("N96").Value = -9494.5
$ws.Range("H98").Value = 1748.75
$ws.Range("I98").Value = 1748.75
$ws.Range("K98").Value = 1748.75
$ws.Range("M98").Value = -250.75
$ws.Range("H101").Value = 14288732
$ws.Range("I101").Value = 25004662
$ws.Range("J101").Value = 825
$ws.Range("K101").Value = 75013986
$ws.Range("L101").Value = 2475
$ws.Range("M101").Value = -75012364
$ws.Range("N101").Value = -5719
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H113").Value = 4483.3335
$ws.Range("I113").Value = 4450
$ws.Range("K113").Value = 4450
$ws.Range("M113").Value = -1196
$ws.Range("H122").Value = 1748.75
$ws.Range("I122").Value = 1748.75
$ws.Range("K122").Value = 5246.25
$ws.Range("M122").Value = -2796.25
$ws.Range("H132").Value = 2856.3333
$ws.Range("I132").Value = 2856.3333
$ws.Range("K132").Value = 8568.999899999999
$ws.Range("M132").Value = -6038.999899999999
$ws.Range("H137").Value = 3120.2942
$ws.Range("I137").Value = 1860.6428
$ws.Range("J137").Value = 8998.666999999999
$ws.Range("K137").Value = 5581.928400000001
$ws.Range("L137").Value = 26996.001
$ws.Range("M137").Value = -3031.928400000001
$ws.Range("N137").Value = -32096.001
$ws.Range("H141").Value = 2008.8572
$ws.Range("I141").Value = 1510.5
$ws.Range("K141").Value = 4531.5
$ws.Range("M141").Value = 648.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3502922
$ws.Range("I32").Value = 3336394.5
$ws.Range("K32").Value = 3336394.5
$ws.Range("M32").Value = -3336107.5
$ws.Range("H61").Value = 3459.6
$ws.Range("I61").Value = 3074.5
$ws.Range("K61").Value = 3074.5
$ws.Range("M61").Value = -2862.5
$ws.Range("H63").Value = 17500.8
$ws.Range("I63").Value = 17500.8
$ws.Range("K63").Value = 17500.8
$ws.Range("M63").Value = -16814.8
$ws.Range("H66").Value = 17500.8
$ws.Range("I66").Value = 17500.8
$ws.Range("K66").Value = 87504
$ws.Range("M66").Value = -84072
$ws.Range("H74").Value = 1482.6923
$ws.Range("I74").Value = 1477.4445
$ws.Range("J74").Value = 1494.5
$ws.Range("K74").Value = 1477.4445
$ws.Range("L74").Value = 1494.5
$ws.Range("M74").Value = -603.4445000000001
$ws.Range("N74").Value = -3242.5
$ws.Range("H77").Value = 1482.6923
$ws.Range("I77").Value = 1477.4445
$ws.Range("J77").Value = 1494.5
$ws.Range("K77").Value = 7387.2225
$ws.Range("L77").Value = 7472.5
$ws.Range("M77").Value = -3019.2225
$ws.Range("N77").Value = -16208.5
$ws.Range("H106").Value = 19999.5
$ws.Range("J106").Value = 19999.5
$ws.Range("L106").Value = 19999.5
$ws.Range("N106").Value = -22523.5
$ws.Range("H122").Value = 4200
$ws.Range("J122").Value = 4200
$ws.Range("L122").Value = 12600
$ws.Range("N122").Value = -17500
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H136").Value = 3459.6
$ws.Range("I136").Value = 3074.5
$ws.Range("K136").Value = 9223.5
$ws.Range("M136").Value = -6673.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3755.3333
$ws.Range("I20").Value = 1971.1428
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 1971.1428
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -1724.1428
$ws.Range("N20").Value = -10494
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H134").Value = 1790
$ws.Range("I134").Value = 1790
$ws.Range("K134").Value = 5370
$ws.Range("M134").Value = -2835

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2497
$ws.Range("I31").Value = 2497
$ws.Range("K31").Value = 2497
$ws.Range("M31").Value = -2202
$ws.Range("H34").Value = 2497
$ws.Range("I34").Value = 2497
$ws.Range("K34").Value = 2497
$ws.Range("M34").Value = -2295
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H105").Value = 3018
$ws.Range("I105").Value = 2156.5
$ws.Range("J105").Value = 4166.6665
$ws.Range("K105").Value = 2156.5
$ws.Range("L105").Value = 4166.6665
$ws.Range("M105").Value = -409.5
$ws.Range("N105").Value = -7660.6665
$ws.Range("H125").Value = 107749.5
$ws.Range("J125").Value = 107749.5
$ws.Range("L125").Value = 107749.5
$ws.Range("N125").Value = -112669.5
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.818184
$ws.Range("I2").Value = 73.5
$ws.Range("K2").Value = 73.5
$ws.Range("M2").Value = 39.5
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H80").Value = 4166.3335
$ws.Range("I80").Value = 3999.5
$ws.Range("K80").Value = 3999.5
$ws.Range("M80").Value = -3001.5
$ws.Range("H83").Value = 4166.3335
$ws.Range("I83").Value = 3999.5
$ws.Range("K83").Value = 19997.5
$ws.Range("M83").Value = -15005.5
$ws.Range("H134").Value = 24598.8
$ws.Range("J134").Value = 24598.8
$ws.Range("L134").Value = 73796.39999999999
$ws.Range("N134").Value = -78866.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8238.357
$ws.Range("J7").Value = 8427.846
$ws.Range("L7").Value = 8427.846
$ws.Range("N7").Value = -8651.846
$ws.Range("H40").Value = 3698.182
$ws.Range("I40").Value = 3011.7144
$ws.Range("K40").Value = 3011.7144
$ws.Range("M40").Value = -2875.7144
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 6234.3794
$ws.Range("I122").Value = 4743.0835
$ws.Range("J122").Value = 7287.0586
$ws.Range("K122").Value = 14229.2505
$ws.Range("L122").Value = 21861.1758
$ws.Range("M122").Value = -11779.2505
$ws.Range("N122").Value = -26761.1758
$ws.Range("H126").Value = 8238.357
$ws.Range("J126").Value = 8427.846
$ws.Range("L126").Value = 25283.538
$ws.Range("N126").Value = -30223.538
$ws.Range("H136").Value = 4400
$ws.Range("I136").Value = 4400
$ws.Range("K136").Value = 13200
$ws.Range("M136").Value = -10650

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5000000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H113").Value = 540.6667
$ws.Range("I113").Value = 475.25
$ws.Range("J113").Value = 671.5
$ws.Range("K113").Value = 1425.75
$ws.Range("L113").Value = 2014.5
$ws.Range("M113").Value = 744.25
$ws.Range("N113").Value = -6354.5
$ws.Range("H136").Value = 1784.3334
$ws.Range("I136").Value = 1226.5
$ws.Range("K136").Value = 3679.5
$ws.Range("M136").Value = -1129.5

